$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 33
$ws.Range("H33").Value = 683.88464
$ws.Range("I33").Value = 353.42856
$ws.Range("J33").Value = 2071.8
$ws.Range("K33").Value = 353.42856
$ws.Range("L33").Value = 2071.8
$ws.Range("M33").Value = -124.42856
$ws.Range("N33").Value = -2529.8

# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents() | Out-Null

# Row 70
$ws.Range("H70").Value = 2802.8333
$ws.Range("I70").Value = 3454.5454
$ws.Range("J70").Value = 2251.3845
$ws.Range("K70").Value = 10363.6362
$ws.Range("L70").Value = 6754.1535
$ws.Range("M70").Value = -10093.6362
$ws.Range("N70").Value = -7294.1535

# Row 73
$ws.Range("H73").Value = 2802.8333
$ws.Range("I73").Value = 3454.5454
$ws.Range("J73").Value = 2251.3845
$ws.Range("K73").Value = 10363.6362
$ws.Range("L73").Value = 6754.1535
$ws.Range("M73").Value = -9427.636200000001
$ws.Range("N73").Value = -8626.1535

# Row 137
$ws.Range("H137").Value = 1682.7693
$ws.Range("I137").Value = 1347.6
$ws.Range("K137").Value = 4042.8
$ws.Range("M137").Value = -1492.8

# Row 138
$ws.Range("H138").Value = 2454.375
$ws.Range("J138").Value = 2205.6667
$ws.Range("L138").Value = 6617.000100000001
$ws.Range("N138").Value = -16897.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 97
$ws.Range("H97").Value = 718
$ws.Range("I97").Value = 484.875
$ws.Range("K97").Value = 484.875
$ws.Range("M97").Value = 11.125

# Row 123
$ws.Range("H123").Value = 10000
$ws.Range("J123").Value = 10000
$ws.Range("L123").Value = 10000
$ws.Range("N123").Value = -19800

# Row 132
$ws.Range("H132").Value = 950
$ws.Range("I132").Value = 950
$ws.Range("K132").Value = 2850
$ws.Range("M132").Value = -320

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 3321
$ws.Range("J20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("N20").Value = -15494

# Row 36
$ws.Range("H36").Value = 2718.75
$ws.Range("I36").Value = 2392.8572
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 2392.8572
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -1858.8572
$ws.Range("N36").Value = -6068

# Row 55
$ws.Range("H55").Value = 84632.664
$ws.Range("J55").Value = 84632.664
$ws.Range("L55").Value = 84632.664
$ws.Range("N55").Value = -85178.664

# Row 75
$ws.Range("H75").Value = 10950
$ws.Range("I75").Value = 10950
$ws.Range("K75").Value = 10950
$ws.Range("M75").Value = -10014

# Row 78
$ws.Range("H78").Value = 10950
$ws.Range("I78").Value = 10950
$ws.Range("K78").Value = 32850
$ws.Range("M78").Value = -28170

# Row 99
$ws.Range("H99").Value = 4703
$ws.Range("I99").Value = 4014.7778
$ws.Range("K99").Value = 4014.7778
$ws.Range("M99").Value = -2516.7778

# Row 134
$ws.Range("H134").Value = 1052.6086
$ws.Range("I134").Value = 962.381
$ws.Range("K134").Value = 2887.143
$ws.Range("M134").Value = -352.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 134
$ws.Range("H134").Value = 2435.5625
$ws.Range("I134").Value = 2097.2727
$ws.Range("J134").Value = 3179.8
$ws.Range("K134").Value = 6291.8181
$ws.Range("L134").Value = 9539.400000000001
$ws.Range("M134").Value = -3756.8181
$ws.Range("N134").Value = -14609.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 70
$ws.Range("H70").Value = 245
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents() | Out-Null

# Row 73
$ws.Range("H73").Value = 245
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents() | Out-Null

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 4639.5
$ws.Range("I70").Value = 4639.5
$ws.Range("K70").Value = 4639.5
$ws.Range("M70").Value = -4369.5

# Row 73
$ws.Range("H73").Value = 4639.5
$ws.Range("I73").Value = 4639.5
$ws.Range("K73").Value = 4639.5
$ws.Range("M73").Value = -3703.5

# Row 80
$ws.Range("H80").Value = 2574.375
$ws.Range("J80").Value = 2686.5
$ws.Range("L80").Value = 2686.5
$ws.Range("N80").Value = -4682.5

# Row 83
$ws.Range("H83").Value = 2574.375
$ws.Range("J83").Value = 2686.5
$ws.Range("L83").Value = 13432.5
$ws.Range("N83").Value = -23416.5

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents() | Out-Null
$ws.Range("N102").ClearContents() | Out-Null

# Row 107
$ws.Range("H107").Value = 511.23077
$ws.Range("I107").Value = 305.3
$ws.Range("J107").Value = 1197.6666
$ws.Range("K107").Value = 305.3
$ws.Range("L107").Value = 1197.6666
$ws.Range("M107").Value = 1614.7
$ws.Range("N107").Value = -5037.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Range("H40").Value = 4967.143
$ws.Range("I40").Value = 4895
$ws.Range("J40").Value = 5400
$ws.Range("K40").Value = 4895
$ws.Range("L40").Value = 5400
$ws.Range("M40").Value = -4759
$ws.Range("N40").Value = -5672

# Row 46
$ws.Range("H46").Value = 1715.5385
$ws.Range("J46").Value = 2250.3333
$ws.Range("L46").Value = 2250.3333
$ws.Range("N46").Value = -2626.3333

# Row 82
$ws.Range("H82").Value = 1093.9
$ws.Range("I82").Value = 802.7143
$ws.Range("J82").Value = 1773.3334
$ws.Range("K82").Value = 802.7143
$ws.Range("L82").Value = 1773.3334
$ws.Range("M82").Value = -441.7143
$ws.Range("N82").Value = -2495.3334

# Row 85
$ws.Range("H85").Value = 1093.9
$ws.Range("I85").Value = 802.7143
$ws.Range("J85").Value = 1773.3334
$ws.Range("K85").Value = 802.7143
$ws.Range("L85").Value = 1773.3334
$ws.Range("M85").Value = 445.2857
$ws.Range("N85").Value = -4269.3334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 62
$ws.Range("H62").Value = 4332.1665
$ws.Range("I62").Value = 4066
$ws.Range("J62").Value = 4598.3335
$ws.Range("K62").Value = 4066
$ws.Range("L62").Value = 4598.3335
$ws.Range("M62").Value = -3442
$ws.Range("N62").Value = -5846.3335

# Row 65
$ws.Range("H65").Value = 4332.1665
$ws.Range("I65").Value = 4066
$ws.Range("J65").Value = 4598.3335
$ws.Range("K65").Value = 20330
$ws.Range("L65").Value = 22991.6675
$ws.Range("M65").Value = -17210
$ws.Range("N65").Value = -29231.6675

# Row 74
$ws.Range("H74").Value = 49034.5
$ws.Range("I74").Value = 63569
$ws.Range("J74").Value = 34500
$ws.Range("K74").Value = 63569
$ws.Range("L74").Value = 34500
$ws.Range("M74").Value = -62633
$ws.Range("N74").Value = -36372

# Row 77
$ws.Range("H77").Value = 49034.5
$ws.Range("I77").Value = 63569
$ws.Range("J77").Value = 34500
$ws.Range("K77").Value = 190707
$ws.Range("L77").Value = 103500
$ws.Range("M77").Value = -186027
$ws.Range("N77").Value = -112860

# Row 107
$ws.Range("H107").Value = 799.6667
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 699.5
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 2098.5
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -5938.5

# Row 126
$ws.Range("H126").Value = 3283.3333
$ws.Range("I126").Value = 4250
$ws.Range("K126").Value = 12750
$ws.Range("M126").Value = -10280

# Row 132
$ws.Range("H132").Value = 2405.1667
$ws.Range("I132").Value = 2252.5293
$ws.Range("K132").Value = 6757.5879
$ws.Range("M132").Value = -4227.5879
